# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# All Price/Volume(1h) cells are stored as plain text in the sheet. Values that
# look like valid numbers (e.g. "206.44") are prefixed here with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cells) instead of silently converting them to numeric values and losing
# formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.926.97'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.555.05'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''206.44'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''21.93'
$ws.Range("E8").Value = '  +2.56%  '
$ws.Range("D9").Value = '''0.247'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +0.75%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '1.775.94'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = '1.554.14'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '26.921.31'
$ws.Range("D17").Value = '''61.71'
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '''218.37'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("E19").Value = '  +1.61%  '
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '''154.10'
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("E30").Value = '  +2.28%  '
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '''3.22'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '1.432.55'
$ws.Range("E33").Value = '  +4.92%  '
$ws.Range("E34").Value = '  +4.39%  '
$ws.Range("D35").Value = '''1.56'
$ws.Range("E35").Value = '  +3.74%  '
$ws.Range("D36").Value = '''0.980'
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = '''0.521'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D45").Value = '''63.85'
$ws.Range("E45").Value = '  +1.10%  '
$ws.Range("D46").Value = '''1.75'
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("D47").Value = '1.690.24'
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").Value = '''86.96'
$ws.Range("D49").Value = '''0.0524'
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("E50").Value = '  +3.80%  '
$ws.Range("E51").Value = '  +1.37%  '
